# Fruta / hortaliza, semanal
# Insert a new weekly record as row 9, shifting the previous rows 9 and 10
# down to rows 10 and 11, then populate the new row 9 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 9 (pushes rows 9-10 to 10-11)
$ws.Rows.Item(9).Insert()

# Fill in the new row 9 with this week's values
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 45062
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112030
$ws.Range("G9").Value = "Poroto granado"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 1700
$ws.Range("K9").Value = 2800
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 2900
$ws.Range("N9").Value = "$/kilo"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 2900
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"
